# Applies the diff to Documentación/CU-01_IniciarSesión/Descripción.docx:
#   1) Adds bookmark "_Hlk178348041" spanning the whole description table
#      (start: right before the "Nombre" run in the table's first cell;
#      end: as the table's own last child, i.e. a sibling of the <w:tr>
#      rows, *after* the final row and *before* </w:tbl>).
#   2) Renames the button label from "Iniciar sesión" to "Entrar" in the
#      three "flujo normal" bullet sentences (the FA-01 error-message
#      sentence keeps saying "Iniciar sesión").

$d = $word.ActiveDocument

# --- 1) Bookmark -------------------------------------------------------
# A <w:bookmarkEnd> that is a direct child of <w:tbl> (not inside any
# paragraph) cannot be produced through Range/Bookmarks.Add, which only
# ever anchors both ends inside a single paragraph. Instead we rebuild
# the table's own OOXML (captured verbatim below, unchanged except for
# the two bookmark markers we splice in) and swap it in with InsertXML,
# which keeps every row/cell byte-for-byte identical to the original.

$tableXml = @'
<w:tbl><w:tblPr><w:tblStyle w:val="Tablaconcuadrcula"/><w:tblW w:w="5000" w:type="pct"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2037"/><w:gridCol w:w="8753"/></w:tblGrid><w:tr w:rsidR="006C0048" w14:paraId="24D0A206" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="6BD79D93" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Nombre</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="49C4FAA7" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>CU-01 Iniciar sesión</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="7D475B9A" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="72E49359" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Descripción</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="45F5AF00" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El caso de uso tiene como finalidad que un actor inicie sesión en el sistema</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="683F4529" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="35922386" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Trazabilidad</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="426999D4" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>FRQ-35, NFRQ-03, NFRQ-08</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="3B49B5B8" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="5E4D2C5B" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Actor(es)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="51AB7F91" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Administrador del supermercado, paquetería, cajero, contador</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="38A16C50" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="15CD17BB" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Disparador</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="78142A23" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El actor inicia la aplicación</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="58C3306E" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="23C50C9C" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Precondiciones</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="59818474" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>PRE-01 El actor debe estar registrado en la base de datos</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="240CDBB0" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="5325C114" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Flujo normal</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="6A6D034A" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El sistema muestra la ventana “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IniciarSesiónView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” con los campos de correo y contraseña y un botón “Iniciar sesión” deshabilitado.</w:t></w:r></w:p><w:p w14:paraId="732424B7" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El actor ingresa su correo y contraseña.</w:t></w:r></w:p><w:p w14:paraId="4CA183D4" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El sistema habilita el botón “Iniciar sesión”.</w:t></w:r></w:p><w:p w14:paraId="6B4B45E3" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El actor da clic en el botón “Iniciar sesión”.</w:t></w:r></w:p><w:p w14:paraId="41C5088E" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El sistema consulta el EMPLEADO en la base de datos (EX-01),</w:t></w:r></w:p><w:p w14:paraId="3EE6096D" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:jc w:val="both"/></w:pPr><w:r><w:t>descifra la contraseña y verifica que las credenciales sean correctas (FA-01),</w:t></w:r></w:p><w:p w14:paraId="1A046C46" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:jc w:val="both"/></w:pPr><w:r><w:t>cierra la ventana “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IniciarSesiónView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” y abre la ventana “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PrincipalView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” (FA-02)</w:t></w:r></w:p><w:p w14:paraId="36BA6C75" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:jc w:val="both"/></w:pPr><w:r><w:t>y agrega el nombre completo del empleado en la ventana de inicio.</w:t></w:r></w:p><w:p w14:paraId="28DF92BB" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Fin del caso de uso.</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="5B0617F4" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="1C54BAEB" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Flujo alterno</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="799383E3" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>FA-01 Datos inválidos</w:t></w:r></w:p><w:p w14:paraId="2B379659" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El sistema muestra el mensaje “Correo y/o contraseña incorrectos, verifique e inténtelo de nuevo” debajo del botón “Iniciar sesión”.</w:t></w:r></w:p><w:p w14:paraId="0E8FA9FA" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Vuelve al flujo normal en el paso 2.</w:t></w:r></w:p><w:p w14:paraId="3878A6E6" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>FA-02 El sistema reconoce al usuario como cajero.</w:t></w:r></w:p><w:p w14:paraId="07A79F8C" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:t>El sistema abre la ventana “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RegistrarVentaView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”.</w:t></w:r></w:p><w:p w14:paraId="44A0771B" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r w:rsidRPr="00204B76"><w:t>Regresa al flujo normal en el último paso.</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="1B69B9F2" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="13D1E0FB" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Excepciones</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="75B0C1E4" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="00987805" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">EX-01 No hay </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">onexión a la </w:t></w:r><w:r><w:t>R</w:t></w:r><w:r w:rsidRPr="00987805"><w:t>ed</w:t></w:r></w:p><w:p w14:paraId="190B3753" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="00987805" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">El sistema muestra la ventana </w:t></w:r><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ErrorView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”</w:t></w:r><w:r w:rsidRPr="00987805"><w:t xml:space="preserve"> con el mensaje “No se pudo conectar a </w:t></w:r><w:r><w:t>la red del supermercado</w:t></w:r><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>inténtelo de nuevo más tarde</w:t></w:r><w:r w:rsidRPr="00987805"><w:t>” junto con un botón de aceptar.</w:t></w:r></w:p><w:p w14:paraId="5E56296B" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="00987805" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">El </w:t></w:r><w:r><w:t xml:space="preserve">cajero </w:t></w:r><w:r w:rsidRPr="00987805"><w:t>da clic en “Aceptar”.</w:t></w:r></w:p><w:p w14:paraId="2A7FC58A" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="00987805" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">El sistema cierra las ventanas </w:t></w:r><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ErrorView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”</w:t></w:r><w:r w:rsidRPr="00987805"><w:t xml:space="preserve"> y </w:t></w:r><w:r><w:t>“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IniciarSesiónView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”</w:t></w:r><w:r w:rsidRPr="00987805"><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p w14:paraId="02CA068B" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="006C0048"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r w:rsidRPr="00987805"><w:t>Regresa al flujo normal en el último paso.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C0048" w14:paraId="3B36A0F9" w14:textId="77777777" w:rsidTr="009C5E11"><w:tc><w:tcPr><w:tcW w:w="944" w:type="pct"/></w:tcPr><w:p w14:paraId="420F3C47" w14:textId="77777777" w:rsidR="006C0048" w:rsidRPr="001F093E" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Postcondiciones</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4056" w:type="pct"/></w:tcPr><w:p w14:paraId="5D884596" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:r><w:t>POS-01 El sistema deja autenticado al EMPLEADO.</w:t></w:r></w:p><w:p w14:paraId="4FE26482" w14:textId="77777777" w:rsidR="006C0048" w:rsidRDefault="006C0048" w:rsidP="009C5E11"><w:r><w:t>POS-02 El sistema autoriza al EMPLEADO a operaciones en base a su PUESTO.</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@

$nombreRun = '<w:r w:rsidRPr="001F093E"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Nombre</w:t></w:r>'
$bookmarkStart = '<w:bookmarkStart w:id="0" w:name="_Hlk178348041"/>'
$bookmarkEnd = '<w:bookmarkEnd w:id="0"/>'

$tableXml = $tableXml.Replace($nombreRun, $bookmarkStart + $nombreRun)

$closeTag = '</w:tr></w:tbl>'
$closeTagPos = $tableXml.LastIndexOf($closeTag)
$tableXml = $tableXml.Substring(0, $closeTagPos) + '</w:tr>' + $bookmarkEnd + '</w:tbl>'

$pkgXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $tableXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$t = $d.Tables.Item(1)
$insertPos = $t.Range.Start
$t.Delete()
$d.Range($insertPos, $insertPos).InsertXML($pkgXml)

# --- 2) Button label text changes --------------------------------------

$d.Content.Find.Execute(
    "” con los campos de correo y contraseña y un botón “Iniciar sesión” deshabilitado.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "” con los campos de correo y contraseña y un botón “Entrar” deshabilitado.",
    2) | Out-Null

$d.Content.Find.Execute(
    "El sistema habilita el botón “Iniciar sesión”.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "El sistema habilita el botón “Entrar”.",
    2) | Out-Null

$d.Content.Find.Execute(
    "El actor da clic en el botón “Iniciar sesión”.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "El actor da clic en el botón “Entrar”.",
    2) | Out-Null

Write-Output "Edit applied."
